# Doing Updates for Financials
# Insert 3 new quarterly columns (D:F) in front of the existing data,
# shifting the old D:K data right to G:N, then fill the new columns
# with the latest quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three blank columns before column D; everything from D:K
#    (dates + financials) shifts right to G:N.
$ws.Columns("D:F").Insert()

# 2) The newly inserted columns come in with column C's formatting
#    (Excel's default Insert behaviour). Re-apply the correct number
#    formats/styles by copying them from the (now shifted) column G,
#    which still holds the original column-D formatting.
$ws.Range("G7:N102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the three "Period Ending" header rows with the three
#    newest quarter-end dates.
$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 4).Value = 43465
    $ws.Cells.Item($r, 5).Value = 43373
    $ws.Cells.Item($r, 6).Value = 43281
}

# 4) "Total Revenue" rows get the new quarter's revenue figures.
$revenueRows = @(8, 18, 32)
foreach ($r in $revenueRows) {
    $ws.Cells.Item($r, 4).Value = 300
    $ws.Cells.Item($r, 5).Value = 900
    $ws.Cells.Item($r, 6).Value = 400
}

# 5) The "Cost of Revenue" row is the negated Total Revenue row.
$negRevenueRows = @(20)
foreach ($r in $negRevenueRows) {
    $ws.Cells.Item($r, 4).Value = -300
    $ws.Cells.Item($r, 5).Value = -900
    $ws.Cells.Item($r, 6).Value = -400
}

# 6) Rows that were "NA" across every quarter stay "NA" in the new
#    columns too.
$naRows = @(9, 10, 12, 21)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}

# 7) Rows that were 0 across every quarter stay 0 in the new columns.
$zeroRows = @(13,14,15,17,22,23,24,25,26,27,28,29,30,31,33,34,35,41,42,43,44,45,46,47,48,49,50,51,52,53,54,57,58,59,60,61,62,63,64,65,66,68,69,70,71,72,73,74,75,76,77,81,83,84,85,86,87,88,89,91,92,93,94,96,97,98,99,100,101,102)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
}

# 8) Section-spacer rows that were blank in every quarter remain blank
#    (nothing to do - Insert already left D:F empty for these rows).
